$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trim trailing space from A2 ("Baz " -> "Baz")
$ws.Range("A2").Value = "Baz"

# Add a thin border around B2 ("Quuk") - creates new border + cellXf
$ws.Range("B2").Borders.LineStyle = 1

# Row 2 height
$ws.Rows.Item(2).RowHeight = 14.9

# Move selection to B2
$ws.Range("B2").Select() | Out-Null
